$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the sales data (B2:M4) with the new random values ---
$row2 = 31, 63, 37, 26, 99, 79, 50, 58, 45, 57, 43, 28
$row3 = 21, 97, 81, 52, 40, 39, 40, 73, 79, 55, 54, 70
$row4 = 61, 56, 8, 72, 20, 25, 91, 96, 9, 10, 64, 79

for ($i = 0; $i -lt 12; $i++) {
    $col = $i + 2
    $ws.Cells.Item(2, $col).Value = $row2[$i]
    $ws.Cells.Item(3, $col).Value = $row3[$i]
    $ws.Cells.Item(4, $col).Value = $row4[$i]
}

# --- Grab the existing line chart on the sheet ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

# --- Chart title ---
$chart.HasTitle = $true
$chart.ChartTitle.Text = "Sales per month"

# --- Axis titles (xlCategory = 1, xlValue = 2, xlPrimary = 1) ---
$catAxis = $chart.Axes(1, 1)
$catAxis.HasTitle = $true
$catAxis.AxisTitle.Text = "Months"

$valAxis = $chart.Axes(2, 1)
$valAxis.HasTitle = $true
$valAxis.AxisTitle.Text = "Sales"

# --- Give every series a category reference pointing at the month headers ---
$catRange = $ws.Range("B1:M1")
$catFormula = "='" + $ws.Name + "'!" + $catRange.Address($true, $true)

$series = $chart.SeriesCollection()
for ($i = 1; $i -le $series.Count; $i++) {
    $s = $series.Item($i)
    $s.XValues = $catFormula
}

Write-Output "chart updated"
